# Update gh-pages to output generated at 456a3b4
# F column holds "想去人数" (want-to-go count) for each event.
# The same events are listed both on their category sheet ("展览")
# and on the aggregate sheet ("全部类型"); bump both occurrences.

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F4").Value = 43
$wsExhibit.Range("F9").Value = 11326
$wsExhibit.Range("F12").Value = 69
$wsExhibit.Range("F14").Value = 5652

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F5").Value = 43
$wsAll.Range("F11").Value = 11326
$wsAll.Range("F14").Value = 69
$wsAll.Range("F17").Value = 5652
